$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value2 = 179.4375
$ws.Cells.Item(33, 9).Value2 = 222.33333
$ws.Cells.Item(33, 10).Value2 = 124.28571
$ws.Cells.Item(33, 11).Value2 = 222.33333
$ws.Cells.Item(33, 12).Value2 = 124.28571
$ws.Cells.Item(33, 13).Value2 = 6.666670000000011
$ws.Cells.Item(33, 14).Value2 = -582.28571

$ws.Cells.Item(43, 8).Value2 = 1000
$ws.Cells.Item(43, 9).Value2 = 1000
$ws.Cells.Item(43, 11).Value2 = 1000
$ws.Cells.Item(43, 13).Value2 = -931

$ws.Cells.Item(69, 8).Value2 = 4275657.5
$ws.Cells.Item(69, 10).Value2 = 8549315
$ws.Cells.Item(69, 12).Value2 = 25647945
$ws.Cells.Item(69, 14).Value2 = -25649693

$ws.Cells.Item(70, 8).Value2 = 1325.3158
$ws.Cells.Item(70, 9).Value2 = 941.63635
$ws.Cells.Item(70, 10).Value2 = 1852.875
$ws.Cells.Item(70, 11).Value2 = 2824.90905
$ws.Cells.Item(70, 12).Value2 = 5558.625
$ws.Cells.Item(70, 13).Value2 = -2554.90905
$ws.Cells.Item(70, 14).Value2 = -6098.625

$ws.Cells.Item(72, 8).Value2 = 4275657.5
$ws.Cells.Item(72, 10).Value2 = 8549315
$ws.Cells.Item(72, 12).Value2 = 76943835
$ws.Cells.Item(72, 14).Value2 = -76952571

$ws.Cells.Item(73, 8).Value2 = 1325.3158
$ws.Cells.Item(73, 9).Value2 = 941.63635
$ws.Cells.Item(73, 10).Value2 = 1852.875
$ws.Cells.Item(73, 11).Value2 = 2824.90905
$ws.Cells.Item(73, 12).Value2 = 5558.625
$ws.Cells.Item(73, 13).Value2 = -1888.90905
$ws.Cells.Item(73, 14).Value2 = -7430.625

$ws.Cells.Item(76, 8).Value2 = 3177677.2
$ws.Cells.Item(76, 9).Value2 = 3971121.8
$ws.Cells.Item(76, 10).Value2 = 3900
$ws.Cells.Item(76, 11).Value2 = 3971121.8
$ws.Cells.Item(76, 12).Value2 = 3900
$ws.Cells.Item(76, 13).Value2 = -3970806.8
$ws.Cells.Item(76, 14).Value2 = -4530

$ws.Cells.Item(79, 8).Value2 = 3177677.2
$ws.Cells.Item(79, 9).Value2 = 3971121.8
$ws.Cells.Item(79, 10).Value2 = 3900
$ws.Cells.Item(79, 11).Value2 = 3971121.8
$ws.Cells.Item(79, 12).Value2 = 3900
$ws.Cells.Item(79, 13).Value2 = -3970029.8
$ws.Cells.Item(79, 14).Value2 = -6084

$ws.Cells.Item(100, 8).Value2 = 16669427
$ws.Cells.Item(100, 9).Value2 = 33334914
$ws.Cells.Item(100, 10).Value2 = 3940
$ws.Cells.Item(100, 11).Value2 = 33334914
$ws.Cells.Item(100, 12).Value2 = 3940
$ws.Cells.Item(100, 13).Value2 = -33334373
$ws.Cells.Item(100, 14).Value2 = -5022

$ws.Cells.Item(111, 8).Value2 = 1850
$ws.Cells.Item(111, 9).Value2 = 700
$ws.Cells.Item(111, 10).Value2 = 3000
$ws.Cells.Item(111, 11).Value2 = 2100
$ws.Cells.Item(111, 12).Value2 = 9000
$ws.Cells.Item(111, 13).Value2 = 967
$ws.Cells.Item(111, 14).Value2 = -15134

$ws.Cells.Item(138, 8).Value2 = 7403795.5
$ws.Cells.Item(138, 9).Value2 = 1787825.8
$ws.Cells.Item(138, 10).Value2 = 12197916
$ws.Cells.Item(138, 11).Value2 = 5363477.4
$ws.Cells.Item(138, 12).Value2 = 36593748
$ws.Cells.Item(138, 13).Value2 = -5358337.4
$ws.Cells.Item(138, 14).Value2 = -36604028

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 17787.607
$ws.Cells.Item(32, 9).Value2 = 3564.324
$ws.Cells.Item(32, 10).Value2 = 73890.55499999999
$ws.Cells.Item(32, 11).Value2 = 3564.324
$ws.Cells.Item(32, 12).Value2 = 73890.55499999999
$ws.Cells.Item(32, 13).Value2 = -3277.324
$ws.Cells.Item(32, 14).Value2 = -74464.55499999999

$ws.Cells.Item(45, 8).Value2 = 995
$ws.Cells.Item(45, 9).Value2 = 995
$ws.Cells.Item(45, 11).Value2 = 995
$ws.Cells.Item(45, 13).Value2 = -618

$ws.Cells.Item(61, 8).Value2 = 2038.0646
$ws.Cells.Item(61, 9).Value2 = 1489.5
$ws.Cells.Item(61, 10).Value2 = 4890.6
$ws.Cells.Item(61, 11).Value2 = 1489.5
$ws.Cells.Item(61, 12).Value2 = 4890.6
$ws.Cells.Item(61, 13).Value2 = -1277.5
$ws.Cells.Item(61, 14).Value2 = -5314.6

$ws.Cells.Item(132, 8).Value2 = 2642.4358
$ws.Cells.Item(132, 9).Value2 = 2188.3333
$ws.Cells.Item(132, 10).Value2 = 5140
$ws.Cells.Item(132, 11).Value2 = 6564.999899999999
$ws.Cells.Item(132, 12).Value2 = 15420
$ws.Cells.Item(132, 13).Value2 = -4034.999899999999
$ws.Cells.Item(132, 14).Value2 = -20480

$ws.Cells.Item(136, 8).Value2 = 2038.0646
$ws.Cells.Item(136, 9).Value2 = 1489.5
$ws.Cells.Item(136, 10).Value2 = 4890.6
$ws.Cells.Item(136, 11).Value2 = 4468.5
$ws.Cells.Item(136, 12).Value2 = 14671.8
$ws.Cells.Item(136, 13).Value2 = -1918.5
$ws.Cells.Item(136, 14).Value2 = -19771.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value2 = 4738.3335
$ws.Cells.Item(82, 9).Value2 = 4738.3335
$ws.Cells.Item(82, 10).Value2 = 0
$ws.Cells.Item(82, 11).Value2 = 4738.3335
$ws.Cells.Item(82, 12).Value2 = 0
$ws.Cells.Item(82, 13).Value2 = -4355.3335
$ws.Cells.Item(82, 14).ClearContents()

$ws.Cells.Item(85, 8).Value2 = 4738.3335
$ws.Cells.Item(85, 9).Value2 = 4738.3335
$ws.Cells.Item(85, 10).Value2 = 0
$ws.Cells.Item(85, 11).Value2 = 4738.3335
$ws.Cells.Item(85, 12).Value2 = 0
$ws.Cells.Item(85, 13).Value2 = -3412.3335
$ws.Cells.Item(85, 14).ClearContents()

$ws.Cells.Item(105, 8).Value2 = 266325.84
$ws.Cells.Item(105, 9).Value2 = 2944.1667
$ws.Cells.Item(105, 10).Value2 = 717837.3
$ws.Cells.Item(105, 11).Value2 = 2944.1667
$ws.Cells.Item(105, 12).Value2 = 717837.3
$ws.Cells.Item(105, 13).Value2 = -1197.1667
$ws.Cells.Item(105, 14).Value2 = -721331.3

$ws.Cells.Item(134, 8).Value2 = 4547.952
$ws.Cells.Item(134, 9).Value2 = 3476.0715
$ws.Cells.Item(134, 10).Value2 = 6691.7144
$ws.Cells.Item(134, 11).Value2 = 10428.2145
$ws.Cells.Item(134, 12).Value2 = 20075.1432
$ws.Cells.Item(134, 13).Value2 = -7893.2145
$ws.Cells.Item(134, 14).Value2 = -25145.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value2 = 790.5263
$ws.Cells.Item(105, 9).Value2 = 790.5263
$ws.Cells.Item(105, 10).Value2 = 0
$ws.Cells.Item(105, 11).Value2 = 790.5263
$ws.Cells.Item(105, 12).Value2 = 0
$ws.Cells.Item(105, 13).Value2 = 956.4737
$ws.Cells.Item(105, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value2 = 21739692
$ws.Cells.Item(113, 9).Value2 = 497
$ws.Cells.Item(113, 10).Value2 = 31250590
$ws.Cells.Item(113, 11).Value2 = 1491
$ws.Cells.Item(113, 12).Value2 = 93751770
$ws.Cells.Item(113, 14).Value2 = -93756110
$ws.Cells.Item(113, 13).Value2 = 679

$ws.Cells.Item(131, 8).Value2 = 1296.6709
$ws.Cells.Item(131, 10).Value2 = 1397.7042
$ws.Cells.Item(131, 12).Value2 = 4193.112599999999
$ws.Cells.Item(131, 14).Value2 = -14273.1126

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value2 = 2505.4688
$ws.Cells.Item(102, 9).Value2 = 2434.3
$ws.Cells.Item(102, 10).Value2 = 2624.0833
$ws.Cells.Item(102, 11).Value2 = 2434.3
$ws.Cells.Item(102, 12).Value2 = 2624.0833
$ws.Cells.Item(102, 13).Value2 = -812.3000000000002
$ws.Cells.Item(102, 14).Value2 = -5868.0833

$ws.Cells.Item(122, 8).Value2 = 1459.909
$ws.Cells.Item(122, 9).Value2 = 1294.3334
$ws.Cells.Item(122, 10).Value2 = 1658.6
$ws.Cells.Item(122, 11).Value2 = 3883.0002
$ws.Cells.Item(122, 12).Value2 = 4975.799999999999
$ws.Cells.Item(122, 13).Value2 = -1433.0002
$ws.Cells.Item(122, 14).Value2 = -9875.799999999999

$ws.Cells.Item(123, 8).Value2 = 11057.565
$ws.Cells.Item(123, 10).Value2 = 11057.565
$ws.Cells.Item(123, 12).Value2 = 11057.565
$ws.Cells.Item(123, 14).Value2 = -15957.565

$ws.Cells.Item(132, 8).Value2 = 3472.5386
$ws.Cells.Item(132, 9).Value2 = 2924.8462
$ws.Cells.Item(132, 10).Value2 = 4567.923
$ws.Cells.Item(132, 11).Value2 = 8774.5386
$ws.Cells.Item(132, 12).Value2 = 13703.769
$ws.Cells.Item(132, 13).Value2 = -6244.5386
$ws.Cells.Item(132, 14).Value2 = -18763.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 864.875
$ws.Cells.Item(46, 9).Value2 = 703.3333
$ws.Cells.Item(46, 10).Value2 = 1349.5
$ws.Cells.Item(46, 11).Value2 = 703.3333
$ws.Cells.Item(46, 12).Value2 = 1349.5
$ws.Cells.Item(46, 13).Value2 = -515.3333
$ws.Cells.Item(46, 14).Value2 = -1725.5

$ws.Cells.Item(122, 8).Value2 = 3490.8262
$ws.Cells.Item(122, 9).Value2 = 1399.5
$ws.Cells.Item(122, 11).Value2 = 4198.5
$ws.Cells.Item(122, 13).Value2 = -1748.5

$ws.Cells.Item(136, 8).Value2 = 2949.375
$ws.Cells.Item(136, 9).Value2 = 1569.7028
$ws.Cells.Item(136, 10).Value2 = 7590.091
$ws.Cells.Item(136, 11).Value2 = 4709.1084
$ws.Cells.Item(136, 12).Value2 = 22770.273
$ws.Cells.Item(136, 13).Value2 = -2159.1084
$ws.Cells.Item(136, 14).Value2 = -27870.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value2 = 14789.655
$ws.Cells.Item(62, 9).Value2 = 14444.5
$ws.Cells.Item(62, 10).Value2 = 15354.454
$ws.Cells.Item(62, 11).Value2 = 14444.5
$ws.Cells.Item(62, 12).Value2 = 15354.454
$ws.Cells.Item(62, 13).Value2 = -13820.5
$ws.Cells.Item(62, 14).Value2 = -16602.454

$ws.Cells.Item(65, 8).Value2 = 14789.655
$ws.Cells.Item(65, 9).Value2 = 14444.5
$ws.Cells.Item(65, 10).Value2 = 15354.454
$ws.Cells.Item(65, 11).Value2 = 72222.5
$ws.Cells.Item(65, 12).Value2 = 76772.27
$ws.Cells.Item(65, 13).Value2 = -69102.5
$ws.Cells.Item(65, 14).Value2 = -83012.27

$ws.Cells.Item(81, 8).Value2 = 3208.4243
$ws.Cells.Item(81, 9).Value2 = 1829.875
$ws.Cells.Item(81, 10).Value2 = 4505.8823
$ws.Cells.Item(81, 11).Value2 = 3659.75
$ws.Cells.Item(81, 12).Value2 = 9011.7646
$ws.Cells.Item(81, 13).Value2 = -2598.75
$ws.Cells.Item(81, 14).Value2 = -11133.7646

$ws.Cells.Item(84, 8).Value2 = 3208.4243
$ws.Cells.Item(84, 9).Value2 = 1829.875
$ws.Cells.Item(84, 10).Value2 = 4505.8823
$ws.Cells.Item(84, 11).Value2 = 18298.75
$ws.Cells.Item(84, 12).Value2 = 45058.823
$ws.Cells.Item(84, 13).Value2 = -12994.75
$ws.Cells.Item(84, 14).Value2 = -55666.823

$ws.Cells.Item(132, 8).Value2 = 2831.1
$ws.Cells.Item(132, 9).Value2 = 2734.6333
$ws.Cells.Item(132, 10).Value2 = 3120.5
$ws.Cells.Item(132, 11).Value2 = 8203.8999
$ws.Cells.Item(132, 12).Value2 = 9361.5
$ws.Cells.Item(132, 13).Value2 = -5673.8999
$ws.Cells.Item(132, 14).Value2 = -14421.5

$ws.Cells.Item(136, 8).Value2 = 1229.7755
$ws.Cells.Item(136, 9).Value2 = 754.2059
$ws.Cells.Item(136, 10).Value2 = 2307.7334
$ws.Cells.Item(136, 11).Value2 = 2262.6177
$ws.Cells.Item(136, 12).Value2 = 6923.2002
$ws.Cells.Item(136, 13).Value2 = 287.3822999999998
$ws.Cells.Item(136, 14).Value2 = -12023.2002
